$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11 ("R40") in the Rule column now reads "1" (kept as text, matching
# the original cell's string type) instead of "R40".
$ws.Range("B11").Value = "'1"
